$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "2025-02-03 09:47:31"
$ws.Range("C3").Value = "2025-02-03 09:47:42"
$ws.Range("C4").Value = "2025-02-03 09:48:29"
$ws.Range("C5").Value = "2025-02-03 09:48:33"
$ws.Range("C6").Value = "2025-02-03 09:48:45"
$ws.Range("C7").Value = "2025-02-03 09:49:27"
$ws.Range("C8").Value = "2025-02-03 09:49:33"
$ws.Range("C9").Value = "2025-02-03 09:49:43"
$ws.Range("C10").Value = "2025-02-03 09:49:50"
$ws.Range("C11").Value = "2025-02-03 09:50:06"
$ws.Range("C12").Value = "2025-02-03 09:50:14"
$ws.Range("C13").Value = "2025-02-03 09:50:16"
$ws.Range("C14").Value = "2025-02-03 09:50:24"
$ws.Range("C15").Value = "2025-02-03 09:50:29"
$ws.Range("C16").Value = "2025-02-03 09:50:35"
$ws.Range("C17").Value = "2025-02-03 09:50:44"
$ws.Range("C18").Value = "2025-02-03 09:51:00"
$ws.Range("C19").Value = "2025-02-03 09:51:09"
$ws.Range("C20").Value = "2025-02-03 09:51:14"
$ws.Range("C21").Value = "2025-02-03 09:51:22"
$ws.Range("C22").Value = "2025-02-03 09:51:27"
$ws.Range("B23").Value = "ERRO"
$ws.Range("C23").Value = "2025-02-03 09:51:40"
$ws.Range("C24").Value = "2025-02-03 09:52:16"
$ws.Range("C25").Value = "2025-02-03 09:52:22"
$ws.Range("C26").Value = "2025-02-03 09:52:36"
$ws.Range("C27").Value = "2025-02-03 09:52:53"
$ws.Range("C28").Value = "2025-02-03 09:53:01"
$ws.Range("C29").Value = "2025-02-03 09:53:08"
$ws.Range("C30").Value = "2025-02-03 09:53:17"
$ws.Range("C31").Value = "2025-02-03 09:53:46"
$ws.Range("C32").Value = "2025-02-03 09:54:00"
$ws.Range("C33").Value = "2025-02-03 09:54:16"
$ws.Range("B34").Value = "SUCESSO"
$ws.Range("C34").Value = "2025-02-03 09:54:35"
$ws.Range("C35").Value = "2025-02-03 09:54:46"
$ws.Range("C36").Value = "2025-02-03 09:54:47"
$ws.Range("C37").Value = "2025-02-03 09:54:57"
$ws.Range("B38").Value = "ERRO"
$ws.Range("C38").Value = "2025-02-03 09:55:09"
$ws.Range("C39").Value = "2025-02-03 09:55:41"
$ws.Range("C40").Value = "2025-02-03 09:55:51"
$ws.Range("C41").Value = "2025-02-03 09:56:04"
$ws.Range("C42").Value = "2025-02-03 09:56:10"
$ws.Range("C43").Value = "2025-02-03 09:56:58"
$ws.Range("C44").Value = "2025-02-03 09:57:09"
$ws.Range("C45").Value = "2025-02-03 09:57:52"
$ws.Range("C46").Value = "2025-02-03 09:58:08"
$ws.Range("C47").Value = "2025-02-03 09:58:23"
